$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44491
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = 800
$ws.Range("N2").Value = "`$/kilo (volumen en unidades)"
$ws.Range("P2").Value = 800

# Row 4
$ws.Range("D4").Value = 44167
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 5000
$ws.Range("N4").Value = "`$/unidad"
$ws.Range("O4").Value = "Región de O'Higgins"
$ws.Range("P4").Value = 5000

# Row 5
$ws.Range("D5").Value = 44167
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 560
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("N5").Value = "`$/unidad"
$ws.Range("O5").Value = "Región de O'Higgins"
$ws.Range("P5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44167
$ws.Range("I6").Value = "Tercera"
$ws.Range("J6").Value = 450
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("N6").Value = "`$/unidad"
$ws.Range("O6").Value = "Región de O'Higgins"
$ws.Range("P6").Value = 2000

# Row 7
$ws.Range("D7").Value = 44497
$ws.Range("J7").Value = 250

# Row 8
$ws.Range("D8").Value = 44504
$ws.Range("J8").Value = 200

# Row 9
$ws.Range("D9").Value = 44305
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2500
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 2500

# Row 10
$ws.Range("D10").Value = 44312
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 180
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2500
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 2500

# Row 11
$ws.Range("D11").Value = 44194
$ws.Range("I11").Value = "Extra"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 3500
$ws.Range("L11").Value = 3500
$ws.Range("M11").Value = 3500
$ws.Range("P11").Value = 3500

# Row 12
$ws.Range("D12").Value = 44194
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("O12").Value = "Región de O'Higgins"
$ws.Range("P12").Value = 3000

# Row 15
$ws.Range("D15").Value = 44223
$ws.Range("H15").Value = "Americana O Klondike"
$ws.Range("I15").Value = "Extra"
$ws.Range("J15").Value = 340
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2500
$ws.Range("N15").Value = "`$/unidad"
$ws.Range("O15").Value = "Región de O'Higgins"
$ws.Range("P15").Value = 2500

# Row 16
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 2000
$ws.Range("P16").Value = 2000

# Row 17
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1500
$ws.Range("P17").Value = 1500

# Row 18
$ws.Range("I18").Value = "Tercera"
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 1000
$ws.Range("P18").Value = 1000

# Row 19
$ws.Range("D19").Value = 44477
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 800
$ws.Range("M19").Value = 800
$ws.Range("N19").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 800

# Row 20
$ws.Range("D20").Value = 44510
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 800
$ws.Range("L20").Value = 800
$ws.Range("M20").Value = 800
$ws.Range("N20").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O20").Value = "Perú"
$ws.Range("P20").Value = 800

# Row 21
$ws.Range("D21").Value = 44495
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 800
$ws.Range("M21").Value = 800
$ws.Range("N21").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 800

# Row 22
$ws.Range("D22").Value = 44483
$ws.Range("J22").Value = 120
